# Applies updated Price (D) and Volume(1h) (E) text values to the
# cryptos worksheet for rows 2-51, matching the commit's refreshed
# scrape snapshot. Values are written as literal text (not re-parsed
# as numbers) so formatting such as trailing zeros, thousand-dot
# separators, and the padded "  -x.xx%  " strings survive exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row entry: spreadsheet row number, new Price text (or $null if
# unchanged), new Volume(1h) text (or $null if unchanged).
$updates = @(
    @{ Row = 2; D = "20.000.63"; E = "  -7.91%  " },
    @{ Row = 3; D = "1.410.25"; E = "  -8.41%  " },
    @{ Row = 4; D = "0.9997"; E = "  -0.14%  " },
    @{ Row = 5; D = "0.9998"; E = "  -0.12%  " },
    @{ Row = 6; D = "273.38"; E = "  -5.63%  " },
    @{ Row = 7; D = "0.3700"; E = "  -5.25%  " },
    @{ Row = 8; D = "0.3070"; E = "  -3.08%  " },
    @{ Row = 9; D = "39.31"; E = "  -8.50%  " },
    @{ Row = 10; D = "0.9996"; E = "  -5.21%  " },
    @{ Row = 11; D = "0.06564"; E = "  -8.50%  " },
    @{ Row = 12; D = "0.9999"; E = "  -0.15%  " },
    @{ Row = 13; D = "5.408"; E = "  -3.70%  " },
    @{ Row = 14; D = "6.175"; E = "  -6.61%  " },
    @{ Row = 15; D = $null; E = "  -8.50%  " },
    @{ Row = 16; D = "1.407.62"; E = "  -8.94%  " },
    @{ Row = 17; D = "0.00001007"; E = "  -8.38%  " },
    @{ Row = 18; D = "0.05766"; E = "  -12.20%  " },
    @{ Row = 19; D = "73.64"; E = "  -11.31%  " },
    @{ Row = 20; D = "0.9994"; E = "  -0.12%  " },
    @{ Row = 21; D = "5.607"; E = "  -8.76%  " },
    @{ Row = 22; D = "14.46"; E = "  -5.72%  " },
    @{ Row = 23; D = "10.85"; E = "  -0.19%  " },
    @{ Row = 24; D = "2.310"; E = "  -4.13%  " },
    @{ Row = 25; D = "19.999.36"; E = "  -7.94%  " },
    @{ Row = 26; D = "2.273"; E = "  -3.57%  " },
    @{ Row = 27; D = "138.63"; E = "  -5.91%  " },
    @{ Row = 28; D = "16.86"; E = "  -8.22%  " },
    @{ Row = 29; D = "1.565.75"; E = "  -9.08%  " },
    @{ Row = 30; D = "108.98"; E = "  -7.18%  " },
    @{ Row = 31; D = "3.825"; E = "  -21.08%  " },
    @{ Row = 32; D = "5.378"; E = "  -8.52%  " },
    @{ Row = 33; D = "0.8545"; E = "  -11.29%  " },
    @{ Row = 34; D = "0.07696"; E = "  -5.82%  " },
    @{ Row = 35; D = "8.444"; E = "  -3.95%  " },
    @{ Row = 36; D = "0.05798"; E = "  -4.36%  " },
    @{ Row = 37; D = "4.806"; E = "  -5.76%  " },
    @{ Row = 38; D = "0.9993"; E = "  -0.14%  " },
    @{ Row = 39; D = "0.1927"; E = "  -5.24%  " },
    @{ Row = 40; D = $null; E = "  -6.58%  " },
    @{ Row = 41; D = "10.30"; E = "  -3.09%  " },
    @{ Row = 42; D = "1.067"; E = "  -9.33%  " },
    @{ Row = 43; D = "1.280"; E = "  -10.46%  " },
    @{ Row = 44; D = "0.5299"; E = "  -7.21%  " },
    @{ Row = 45; D = "3.532"; E = "  -5.44%  " },
    @{ Row = 46; D = "12.14"; E = "  -5.90%  " },
    @{ Row = 47; D = "0.5121"; E = "  -6.38%  " },
    @{ Row = 48; D = "1.806"; E = "  -2.95%  " },
    @{ Row = 49; D = "110.17"; E = "  -5.03%  " },
    @{ Row = 50; D = "1.045"; E = "  -9.99%  " },
    @{ Row = 51; D = "0.9996"; E = "  -0.14%  " }
)

foreach ($item in $updates) {
    if ($null -ne $item.D) {
        $cell = $ws.Range("D" + $item.Row)
        # Force text interpretation so strings like "0.3700" or
        # "20.000.63" are not coerced into numbers and lose their
        # exact printed form.
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        # Drop the temporary text format again so the cell keeps
        # the workbook default (General) style, as in the original.
        $cell.ClearFormats()
    }
    if ($null -ne $item.E) {
        $cell = $ws.Range("E" + $item.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $item.E
        $cell.ClearFormats()
    }
}
